# feat: add 2022-Q1 data
#
# - Inserts a new "2022-Q1" worksheet (per-fund holding detail) right
#   before the "总计" summary sheet.
# - Inserts a matching summary row at the top of "总计"'s data table.

$wb = $excel.ActiveWorkbook

# A sheet that already has the "header row" / "index column" formatting
# we want to reuse (bold font, thin box border, centered) so new cells
# pick up the very same style entries instead of inventing new ones.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet right before "总计"
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$ws = $wb.Worksheets.Add($total)
$ws.Name = "2022-Q1"

# Header row values
$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

# Fund rows: code, name, size, stock position, position ratio, market value, rank
$rows = @(
    @("166019", "中欧价值智选回报混合A", "156.17", "94.14", "3.90", "6.0906", 6),
    @("013220", "中欧新兴价值一年持有混合A", "63.33", "94.47", "3.30", "2.0899", 7),
    @("004235", "中欧价值智选回报混合C", "36.40", "94.14", "3.90", "1.4196", 6),
    @("004814", "中欧红利优享灵活配置混合A", "22.96", "93.96", "4.27", "0.9804", 3),
    @("001887", "中欧价值智选回报混合E", "20.77", "94.14", "3.90", "0.8100", 6),
    @("001556", "天弘中证500指数增强A", "41.41", "94.29", "1.69", "0.6998", 9),
    @("013221", "中欧新兴价值一年持有混合C", "16.75", "94.47", "3.30", "0.5528", 7),
    @("004848", "中欧睿泓定期开放灵活配置混合", "23.30", "59.08", "2.04", "0.4753", 10),
    @("004815", "中欧红利优享灵活配置混合C", "8.45", "93.96", "4.27", "0.3608", 3),
    @("006682", "景顺长城中证500指数增强", "16.63", "87.75", "2.04", "0.3393", 4),
    @("001557", "天弘中证500指数增强C", "13.97", "94.29", "1.69", "0.2361", 9),
    @("000978", "景顺长城量化精选股票", "8.51", "93.86", "2.06", "0.1753", 2),
    @("008851", "景顺长城量化对冲策略三个月定期开放灵活配置混合", "5.05", "74.55", "1.62", "0.0818", 4),
    @("012878", "中信建投量化精选6个月持有期混合型证券投资基金A", "6.32", "88.35", "0.85", "0.0537", 10),
    @("012879", "中信建投量化精选6个月持有期混合型证券投资基金C", "4.26", "88.35", "0.85", "0.0362", 10),
    @("009608", "广发中证500指数增强A", "1.75", "93.00", "1.12", "0.0196", 8),
    @("590007", "中邮中证500指数增强A", "0.43", "91.51", "1.39", "0.0060", 9),
    @("009609", "广发中证500指数增强C", "0.43", "93.00", "1.12", "0.0048", 8),
    @("008124", "中邮中证500指数增强C", "0.04", "91.51", "1.39", "0.0006", 9)
)

$lastRow = 1 + $rows.Count
$r = 2
foreach ($row in $rows) {
    $ws.Range("A$r").Value = ($r - 2)
    # Leading apostrophe forces text storage for numeric-looking values
    # (fund code / size / position% / market value), matching the
    # original "number kept as text" data export.
    $ws.Range("B$r").Value = "'" + $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = "'" + $row[2]
    $ws.Range("E$r").Value = "'" + $row[3]
    $ws.Range("F$r").Value = "'" + $row[4]
    $ws.Range("G$r").Value = "'" + $row[5]
    $ws.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# Reuse the existing header / index-column styles (bold, boxed, centered)
# instead of leaving the new cells unstyled.
$styleSrc.Range("B1:G1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$styleSrc.Range("B1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$styleSrc.Range("A2").Copy()
$ws.Range("A2:A$lastRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Drop the "quote-prefix" styling picked up from the leading-apostrophe
# text entry above, so the text cells stay on the default (unstyled) xf
# like the rest of the data columns.
$ws.Range("B2:G$lastRow").Style = "Normal"

# ---------------------------------------------------------------------
# 2. Insert a new summary row for "2022-Q1" at the top of "总计"'s data
#    (re-fetch "总计" by name since inserting a sheet before it shifts
#    its position, and sheet handles here resolve by position)
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 19
$total.Range("D2").Value = 14.43

# Row-insert copies the header row's bold/boxed formatting onto the new
# row by default; B2:D2 should stay unstyled like the rest of the data
# rows, while A2 (the index column) should pick up the bold/boxed style.
$total.Range("B2:D2").Style = "Normal"
$styleSrc.Range("A2").Copy()
$total.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Renumber the row-index column (A) for the rows pushed down, so it stays
# a contiguous 0-based sequence (matches the source data pipeline's output)
for ($r = 3; $r -le 7; $r++) {
    $total.Range("A$r").Value = $r - 2
}
